# Auto-generated Excel COM-interop script
# Commit: "issue #5: property land done"
#
# 1) Cleans stray whitespace / punctuation noise that had crept into the
#    Chinese shared-string text across every sheet (e.g. "0670-0000" ->
#    "06700000", "李桐豪 " variants, "96年08月 20日" -> "96年08月20日", ...).
# 2) Adds 7 new metadata columns to the "土地" (land) sheet: property_category,
#    category, date, legislator_name, legislator_id, source_file, index -
#    mirroring the scraper row metadata already used elsewhere in the pipeline.

$wb = $excel.ActiveWorkbook

# ---- 土地 (sheet1) ----
$ws = $wb.Worksheets.Item("土地")

# Stamp the 7 new columns (I:O) with the same header/data styling
# already used by the existing A:H columns before filling them in.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:O1").PasteSpecial(-4122) | Out-Null
$ws.Range("H2").Copy() | Out-Null
$ws.Range("I2:O5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "area"
$ws.Range("D1").Value = "share_portion"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "acquire_value"
$ws.Range("I1").Value = "property_category"
$ws.Range("J1").Value = "category"
$ws.Range("K1").Value = "date"
$ws.Range("L1").Value = "legislator_name"
$ws.Range("M1").Value = "legislator_id"
$ws.Range("N1").Value = "source_file"
$ws.Range("O1").Value = "index"

$ws.Range("A2").Value = 14
$ws.Range("B2").Value = "臺北市文山區萬芳段三小段06700000地號"
$ws.Range("C2").Value = 237
$ws.Range("D2").Value = "8分之1"
$ws.Range("E2").Value = "李桐豪"
$ws.Range("F2").Value = "83年07月04日"
$ws.Range("G2").Value = "買賣"
$ws.Range("H2").Value = "(超過五年）"
$ws.Range("I2").Value = "land"
$ws.Range("J2").Value = "normal"
$ws.Range("K2").Value = "2013-12-26"
$ws.Range("L2").Value = "李桐豪"
$ws.Range("M2").Value = 896
$ws.Range("N2").Value = "tmp2e9d1"
$ws.Range("O2").Value = 14
$ws.Range("A3").Value = 15
$ws.Range("B3").Value = "新北市新店區華城二段01220003地號"
$ws.Range("C3").Value = 1136.93
$ws.Range("D3").Value = "70000分之14230"
$ws.Range("E3").Value = "黃素香"
$ws.Range("F3").Value = "96年08月20日"
$ws.Range("G3").Value = "買賣"
$ws.Range("H3").Value = "(超過五年）"
$ws.Range("I3").Value = "land"
$ws.Range("J3").Value = "normal"
$ws.Range("K3").Value = "2013-12-26"
$ws.Range("L3").Value = "李桐豪"
$ws.Range("M3").Value = 896
$ws.Range("N3").Value = "tmp2e9d1"
$ws.Range("O3").Value = 15
$ws.Range("A4").Value = 16
$ws.Range("B4").Value = "新北市新店區環河段04070000地號"
$ws.Range("C4").Value = 11362.89
$ws.Range("D4").Value = "100000分之55"
$ws.Range("E4").Value = "李桐豪"
$ws.Range("F4").Value = "102年03月11曰"
$ws.Range("G4").Value = "買賣"
$ws.Range("H4").Value = 571837
$ws.Range("I4").Value = "land"
$ws.Range("J4").Value = "normal"
$ws.Range("K4").Value = "2013-12-26"
$ws.Range("L4").Value = "李桐豪"
$ws.Range("M4").Value = 896
$ws.Range("N4").Value = "tmp2e9d1"
$ws.Range("O4").Value = 16
$ws.Range("A5").Value = 17
$ws.Range("B5").Value = "新北市新店區環河段04070001地號"
$ws.Range("C5").Value = 6.82
$ws.Range("D5").Value = "100000分之55"
$ws.Range("E5").Value = "李桐豪"
$ws.Range("F5").Value = "102年03月11曰"
$ws.Range("G5").Value = "買賣"
$ws.Range("H5").Value = 343
$ws.Range("I5").Value = "land"
$ws.Range("J5").Value = "normal"
$ws.Range("K5").Value = "2013-12-26"
$ws.Range("L5").Value = "李桐豪"
$ws.Range("M5").Value = 896
$ws.Range("N5").Value = "tmp2e9d1"
$ws.Range("O5").Value = 17

# ---- 建物 (sheet2) ----
$ws = $wb.Worksheets.Item("建物")

$ws.Range("B1").Value = "建物標示"
$ws.Range("C1").Value = "area"
$ws.Range("D1").Value = "share_portion"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "acquire_value"

$ws.Range("A2").Value = 22
$ws.Range("B2").Value = "臺北市文山區萬芳段三小段02068000建號"
$ws.Range("C2").Value = 75.73
$ws.Range("D2").Value = "全部"
$ws.Range("E2").Value = "李桐豪"
$ws.Range("F2").Value = "83年07月04日"
$ws.Range("G2").Value = "買賣"
$ws.Range("H2").Value = "(超過五年)"
$ws.Range("A3").Value = 23
$ws.Range("B3").Value = "臺北市文山區萬芳段三小段02073000建號"
$ws.Range("C3").Value = 47.8
$ws.Range("D3").Value = "8分之1"
$ws.Range("E3").Value = "李桐豪"
$ws.Range("F3").Value = "83年07月04日"
$ws.Range("G3").Value = "買賣"
$ws.Range("H3").Value = "(超過五年）"
$ws.Range("A4").Value = 24
$ws.Range("B4").Value = "新北市新店區華城二段00562000建號"
$ws.Range("C4").Value = 121.16
$ws.Range("D4").Value = "全部"
$ws.Range("E4").Value = "黃素香"
$ws.Range("F4").Value = "96年08月20日"
$ws.Range("G4").Value = "買賣"
$ws.Range("H4").Value = "(超過五年）"
$ws.Range("A5").Value = 25
$ws.Range("B5").Value = "新北市新店區環河段00970000建號"
$ws.Range("C5").Value = 57.21
$ws.Range("D5").Value = "全部"
$ws.Range("E5").Value = "李桐豪"
$ws.Range("F5").Value = "102年03月11曰"
$ws.Range("G5").Value = "買賣"
$ws.Range("H5").Value = "12177820(陽台9.24平方公尺）"

# ---- 汽車 (sheet3) ----
$ws = $wb.Worksheets.Item("汽車")

$ws.Range("B1").Value = "廠牌型號"
$ws.Range("C1").Value = "汽缸容量"
$ws.Range("D1").Value = "所有人"
$ws.Range("E1").Value = "登記（取得）時間"
$ws.Range("F1").Value = "登記（取得）原因"
$ws.Range("G1").Value = "取得價額"

$ws.Range("A2").Value = 35
$ws.Range("B2").Value = "中華"
$ws.Range("C2").Value = 2000
$ws.Range("D2").Value = "李桐豪"
$ws.Range("E2").Value = "96年11月20日"
$ws.Range("F2").Value = "買賣"
$ws.Range("G2").Value = 605000
$ws.Range("A3").Value = 36
$ws.Range("B3").Value = "福特"
$ws.Range("C3").Value = 1600
$ws.Range("D3").Value = "黃素香"
$ws.Range("E3").Value = "102年03月26曰"
$ws.Range("F3").Value = "買賣"
$ws.Range("G3").Value = 500000

# ---- 存款 (sheet4) ----
$ws = $wb.Worksheets.Item("存款")

$ws.Range("B1").Value = "存放機構(應敘明分支機構）"
$ws.Range("C1").Value = "種類"
$ws.Range("D1").Value = "幣別"
$ws.Range("E1").Value = "所有人"
$ws.Range("F1").Value = "新臺幣總額.或折合新臺幣總額"

$ws.Range("A2").Value = 51
$ws.Range("B2").Value = "新店中正郵局(第52支局）"
$ws.Range("C2").Value = "活期儲蓄存款"
$ws.Range("D2").Value = "新臺幣"
$ws.Range("E2").Value = "李桐豪"
$ws.Range("F2").Value = 409902
$ws.Range("A3").Value = 52
$ws.Range("B3").Value = "曰盛國際商業銀行內湖分行"
$ws.Range("C3").Value = "活期儲蓄存款"
$ws.Range("D3").Value = "新臺幣"
$ws.Range("E3").Value = "李桐豪"
$ws.Range("F3").Value = 373426
$ws.Range("A4").Value = 53
$ws.Range("B4").Value = "第一商業銀行木柵分行"
$ws.Range("C4").Value = "活期儲蓄存款."
$ws.Range("D4").Value = "新臺幣"
$ws.Range("E4").Value = "李桐豪"
$ws.Range("F4").Value = 520229
$ws.Range("A5").Value = 54
$ws.Range("B5").Value = "台北富邦商業銀行木柵分行"
$ws.Range("C5").Value = "活期儲蓄存款"
$ws.Range("D5").Value = "新臺幣"
$ws.Range("E5").Value = "李桐豪"
$ws.Range("F5").Value = 367001
$ws.Range("A6").Value = 55
$ws.Range("B6").Value = "臺灣銀行木柵分行"
$ws.Range("C6").Value = "活期儲蓄存款"
$ws.Range("D6").Value = "新臺幣"
$ws.Range("E6").Value = "李桐豪"
$ws.Range("F6").Value = 843898
$ws.Range("A7").Value = 56
$ws.Range("B7").Value = "合作金庫商業銀行西門分行"
$ws.Range("C7").Value = "活期儲蓄存款"
$ws.Range("D7").Value = "新臺幣"
$ws.Range("E7").Value = "李桐豪"
$ws.Range("F7").Value = 1325

# ---- 保險 (sheet5) ----
$ws = $wb.Worksheets.Item("保險")

$ws.Range("B1").Value = "保險公司"
$ws.Range("C1").Value = "保險名稱"
$ws.Range("D1").Value = "要保人"
$ws.Range("E1").Value = "備註"

$ws.Range("A2").Value = 91
$ws.Range("B2").Value = "台灣人壽"
$ws.Range("C2").Value = "健康安心終身醫療B型"
$ws.Range("D2").Value = "李桐豪"
$ws.Range("E2").Value = "共2筆"
$ws.Range("A3").Value = 92
$ws.Range("B3").Value = "台灣人壽"
$ws.Range("C3").Value = "歲歲長泰還本終身"
$ws.Range("D3").Value = "李桐豪"
$ws.Range("A4").Value = 93
$ws.Range("B4").Value = "台灣人壽"
$ws.Range("C4").Value = "六六大順增額終身"
$ws.Range("D4").Value = "黃素香"
$ws.Range("A5").Value = 94
$ws.Range("B5").Value = "台灣人壽"
$ws.Range("C5").Value = "金寶貝兒童终身"
$ws.Range("D5").Value = "黃素香"
$ws.Range("A6").Value = 95
$ws.Range("B6").Value = "台灣人壽"
$ws.Range("C6").Value = "喜福還本定期"
$ws.Range("D6").Value = "黃素香"
$ws.Range("E6").Value = "共2筆"
$ws.Range("A7").Value = 96
$ws.Range("B7").Value = "台灣人壽"
$ws.Range("C7").Value = "金如意還本終身"
$ws.Range("D7").Value = "黃素香"
$ws.Range("E7").Value = "共2筆"

# ---- 債務 (sheet6) ----
$ws = $wb.Worksheets.Item("債務")

$ws.Range("B1").Value = "種類"
$ws.Range("C1").Value = "債務人"
$ws.Range("D1").Value = "債權人及地址"
$ws.Range("E1").Value = "餘額"
$ws.Range("F1").Value = "取得（發生）時間"
$ws.Range("G1").Value = "取得（發生）原因"

$ws.Range("A2").Value = 106
$ws.Range("B2").Value = "抵押"
$ws.Range("C2").Value = "李桐豪"
$ws.Range("D2").Value = "台北富邦銀行臺北市中山區中山北路二段50號"
$ws.Range("E2").Value = 3322620
$ws.Range("F2").Value = "102年03月11曰"
$ws.Range("G2").Value = "購買房舍貸款"
